$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 868, shifting the
# existing data (old rows 868:987) down to 870:989.
$ws.Rows("868:869").Insert()

# --- New row 868: "Primera" quality entry dated 45077 (2023-05-31) ---
$ws.Range("A868").Value = 3
$ws.Range("B868").Value = "Femacal de La Calera"
$ws.Range("C868").Value = "Coquimbo"
$ws.Range("D868").Value = 45077
$ws.Range("E868").Value = 5
$ws.Range("F868").Value = 100112006
$ws.Range("G868").Value = "Repollo"
$ws.Range("H868").Value = "Crespo record"
$ws.Range("I868").Value = "Primera"
$ws.Range("J868").Value = 2700
$ws.Range("K868").Value = 1000
$ws.Range("L868").Value = 1100
$ws.Range("M868").Value = 1056
$ws.Range("N868").Value = "$/unidad"
$ws.Range("O868").Value = "Provincia de Quillota"
$ws.Range("P868").Value = 1056
$ws.Range("Q868").Value = 1
$ws.Range("R868").Value = "Hortaliza"

# --- New row 869: "Segunda" quality entry dated 45077 (2023-05-31) ---
$ws.Range("A869").Value = 3
$ws.Range("B869").Value = "Femacal de La Calera"
$ws.Range("C869").Value = "Coquimbo"
$ws.Range("D869").Value = 45077
$ws.Range("E869").Value = 5
$ws.Range("F869").Value = 100112006
$ws.Range("G869").Value = "Repollo"
$ws.Range("H869").Value = "Crespo record"
$ws.Range("I869").Value = "Segunda"
$ws.Range("J869").Value = 1600
$ws.Range("K869").Value = 900
$ws.Range("L869").Value = 900
$ws.Range("M869").Value = 900
$ws.Range("N869").Value = "$/unidad"
$ws.Range("O869").Value = "Provincia de Quillota"
$ws.Range("P869").Value = 900
$ws.Range("Q869").Value = 1
$ws.Range("R869").Value = "Hortaliza"

# Match the date-formatted number format used by the rest of column D.
$ws.Range("D868").NumberFormat = $ws.Range("D870").NumberFormat
$ws.Range("D869").NumberFormat = $ws.Range("D870").NumberFormat
